$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 461.7
$ws.Range("I28").Value = 333.29413
$ws.Range("K28").Value = 333.29413
$ws.Range("M28").Value = 151.70587
# Row 75: Tomes Roam on the Range | Dhalmelskin Codex
$ws.Range("H75").Value = 33024.11
$ws.Range("J75").Value = 33024.11
$ws.Range("L75").Value = 33024.11
$ws.Range("N75").Value = -34896.11
# Row 78: Field Trip to the Unknown (L) | Dhalmelskin Codex
$ws.Range("H78").Value = 33024.11
$ws.Range("J78").Value = 33024.11
$ws.Range("L78").Value = 99072.33
$ws.Range("N78").Value = -108432.33
# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 3386.125
$ws.Range("I86").Value = 2611
$ws.Range("J86").Value = 3773.6875
$ws.Range("K86").Value = 2611
$ws.Range("L86").Value = 3773.6875
$ws.Range("M86").Value = -1488
$ws.Range("N86").Value = -6019.6875
# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 3386.125
$ws.Range("I89").Value = 2611
$ws.Range("J89").Value = 3773.6875
$ws.Range("K89").Value = 13055
$ws.Range("L89").Value = 18868.4375
$ws.Range("M89").Value = -7439
$ws.Range("N89").Value = -30100.4375
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1555.2858
$ws.Range("I132").Value = 1344.7894
$ws.Range("K132").Value = 4034.3682
$ws.Range("M132").Value = -1504.3682
# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 1461.1428
$ws.Range("I135").Value = 871.3333
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 7841.9997
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -5306.9997
$ws.Range("N135").Value = -50070
# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 1577.159
$ws.Range("J137").Value = 2542.2
$ws.Range("L137").Value = 7626.599999999999
$ws.Range("N137").Value = -12726.6
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2814.9792
$ws.Range("I138").Value = 741
$ws.Range("J138").Value = 4173.793
$ws.Range("K138").Value = 2223
$ws.Range("L138").Value = 12521.379
$ws.Range("M138").Value = 2917
$ws.Range("N138").Value = -22801.379
# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 110387.336
$ws.Range("I141").Value = 123435.75
$ws.Range("K141").Value = 370307.25
$ws.Range("M141").Value = -365127.25

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate | Bronze Plate
$ws.Range("H4").Value = 562.5714
$ws.Range("I4").Value = 389.66666
$ws.Range("K4").Value = 389.66666
$ws.Range("M4").Value = -273.66666
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 3131327.5
$ws.Range("I32").Value = 3232048
$ws.Range("K32").Value = 3232048
$ws.Range("M32").Value = -3231761
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1470.7894
$ws.Range("I61").Value = 1289.4706
$ws.Range("K61").Value = 1289.4706
$ws.Range("M61").Value = -1077.4706
# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1947.8445
$ws.Range("I74").Value = 1235.3077
$ws.Range("K74").Value = 1235.3077
$ws.Range("M74").Value = -361.3077000000001
# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1947.8445
$ws.Range("I77").Value = 1235.3077
$ws.Range("K77").Value = 6176.538500000001
$ws.Range("M77").Value = -1808.538500000001
# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2604.5881
$ws.Range("I132").Value = 2023.2142
$ws.Range("K132").Value = 6069.642599999999
$ws.Range("M132").Value = -3539.642599999999
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1470.7894
$ws.Range("I136").Value = 1289.4706
$ws.Range("K136").Value = 3868.4118
$ws.Range("M136").Value = -1318.4118

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 929.38464
$ws.Range("I99").Value = 881.8333
$ws.Range("K99").Value = 881.8333
$ws.Range("M99").Value = 616.1667

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1638.5333
$ws.Range("I58").Value = 1317.85
$ws.Range("K58").Value = 1317.85
$ws.Range("M58").Value = -1114.85
# Row 92: Walk the Walk | Beech Rod
$ws.Range("H92").Value = 13998.333
$ws.Range("J92").Value = 13998.333
$ws.Range("L92").Value = 13998.333
$ws.Range("N92").Value = -18990.333
# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 4969.3335
$ws.Range("I122").Value = 5251.9
$ws.Range("J122").Value = 3556.5
$ws.Range("K122").Value = 15755.7
$ws.Range("L122").Value = 10669.5
$ws.Range("M122").Value = -13305.7
$ws.Range("N122").Value = -15569.5
# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1638.5333
$ws.Range("I136").Value = 1317.85
$ws.Range("K136").Value = 3953.55
$ws.Range("M136").Value = -1403.55

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 1480.875
$ws.Range("I5").Value = 1141.1666
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 3423.4998
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -3311.4998
$ws.Range("N5").Value = -7724
# Row 44: No More Dumpster Diving | Knight's Bread
$ws.Range("H44").Value = 922.5
$ws.Range("I44").Value = 922.5
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 2767.5
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -2369.5
$ws.Range("N44").ClearContents()
# Row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 3483.45
$ws.Range("I55").Value = 646.4
$ws.Range("J55").Value = 4429.1333
$ws.Range("K55").Value = 1939.2
$ws.Range("L55").Value = 13287.3999
$ws.Range("M55").Value = -1762.2
$ws.Range("N55").Value = -13641.3999
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 16920.615
$ws.Range("J131").Value = 1597.4333
$ws.Range("L131").Value = 4792.2999
$ws.Range("N131").Value = -14872.2999
# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 9055.333000000001
$ws.Range("I132").Value = 2935.5715
$ws.Range("J132").Value = 30474.5
$ws.Range("K132").Value = 26420.1435
$ws.Range("L132").Value = 274270.5
$ws.Range("M132").Value = -23890.1435
$ws.Range("N132").Value = -279330.5
# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 1480.875
$ws.Range("I135").Value = 1141.1666
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 10270.4994
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -7735.499400000001
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 5698.154
$ws.Range("I80").Value = 2005
$ws.Range("J80").Value = 6005.9165
$ws.Range("K80").Value = 2005
$ws.Range("L80").Value = 6005.9165
$ws.Range("M80").Value = -1007
$ws.Range("N80").Value = -8001.9165
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 5698.154
$ws.Range("I83").Value = 2005
$ws.Range("J83").Value = 6005.9165
$ws.Range("K83").Value = 10025
$ws.Range("L83").Value = 30029.5825
$ws.Range("M83").Value = -5033
$ws.Range("N83").Value = -40013.5825
# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 2716.9
$ws.Range("I97").Value = 2967.111
$ws.Range("K97").Value = 2967.111
$ws.Range("M97").Value = -2471.111
# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 2497.4614
$ws.Range("I113").Value = 2385.5
$ws.Range("K113").Value = 2385.5
$ws.Range("M113").Value = -215.5

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 307.07144
$ws.Range("I55").Value = 332.72726
$ws.Range("J55").Value = 213
$ws.Range("K55").Value = 332.72726
$ws.Range("L55").Value = 213
$ws.Range("M55").Value = -159.72726
$ws.Range("N55").Value = -559
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 5042
$ws.Range("I122").Value = 3505.1667
$ws.Range("K122").Value = 10515.5001
$ws.Range("M122").Value = -8065.500100000001
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 5485.6665
$ws.Range("I132").Value = 3108.5833
$ws.Range("J132").Value = 14994
$ws.Range("K132").Value = 9325.749899999999
$ws.Range("L132").Value = 44982
$ws.Range("M132").Value = -6795.749899999999
$ws.Range("N132").Value = -50042
# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("I136").Value = 3730.3462
$ws.Range("K136").Value = 11191.0386
$ws.Range("M136").Value = -8641.0386

$ws = $wb.Worksheets.Item("WVR")
# Row 63: Protecting the Foundation | Rainbow Slops of Aiming
$ws.Range("H63").Value = 56662.332
$ws.Range("J63").Value = 59993.5
$ws.Range("L63").Value = 59993.5
$ws.Range("N63").Value = -61241.5
# Row 66: Curb the Gnawing Feeling (L) | Rainbow Slops of Aiming
$ws.Range("H66").Value = 56662.332
$ws.Range("J66").Value = 59993.5
$ws.Range("L66").Value = 179980.5
$ws.Range("N66").Value = -186220.5
# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 1632.5
$ws.Range("I113").Value = 828.9375
$ws.Range("K113").Value = 2486.8125
$ws.Range("M113").Value = -316.8125
# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 3923.5557
$ws.Range("I122").Value = 3923.5557
$ws.Range("K122").Value = 11770.6671
$ws.Range("M122").Value = -9320.667099999999
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 2730.7058
$ws.Range("I126").Value = 2922.0908
$ws.Range("J126").Value = 2379.8333
$ws.Range("K126").Value = 8766.2724
$ws.Range("L126").Value = 7139.499899999999
$ws.Range("M126").Value = -6296.2724
$ws.Range("N126").Value = -12079.4999
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2653.1785
$ws.Range("I132").Value = 2653.1785
$ws.Range("K132").Value = 7959.5355
$ws.Range("M132").Value = -5429.5355
